$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: handback status text changed for both locales
# ("In Translation" -> "Handed back: in sync with en-US"). This is a shared
# string used by E2/F2/E3/F3, so updating all four keeps a single shared
# string (matches the source diff which only edits the <si> text).
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# zh-cn sheet: record the handback target file / handback file / handback
# datetime for both rows, and hyperlink the new "Latest Target File" cells.
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("I2").Value = "a0644925-81bd-4555-b74d-13847df737e1.md"
$zhcn.Range("J2").Value = "a0644925-81bd-4555-b74d-13847df737e1.01015c1d51e9dabb5a2674fb5c3e5b1e701e5437.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-29 00:25:11"

$zhcn.Range("I3").Value = "ed37b0c1-f71a-45ac-954a-94a9a70c5282.md"
$zhcn.Range("J3").Value = "ed37b0c1-f71a-45ac-954a-94a9a70c5282.d23391ddd3bd083b0c0eccc8d10aeaba46d43b77.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-29 00:25:11"

# Recreate the hyperlinks in row-major order (A2, I2, A3, I3) so the new
# "Latest Target File" links for row 2 and row 3 point at the handoff doc,
# same as the existing "Source File Name" links.
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e9ddad79c30fb16c37f9428b0a1306f903cea0f6/e2e/a0644925-81bd-4555-b74d-13847df737e1.md", "", "", "a0644925-81bd-4555-b74d-13847df737e1.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e9ddad79c30fb16c37f9428b0a1306f903cea0f6/e2e/a0644925-81bd-4555-b74d-13847df737e1.md", "", "", "a0644925-81bd-4555-b74d-13847df737e1.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e9ddad79c30fb16c37f9428b0a1306f903cea0f6/e2e/ed37b0c1-f71a-45ac-954a-94a9a70c5282.md", "", "", "ed37b0c1-f71a-45ac-954a-94a9a70c5282.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e9ddad79c30fb16c37f9428b0a1306f903cea0f6/e2e/ed37b0c1-f71a-45ac-954a-94a9a70c5282.md", "", "", "ed37b0c1-f71a-45ac-954a-94a9a70c5282.md")

# Widen the columns that now hold long file names / hyperlinked text
# (ColumnWidth is character-width based and the host re-quantizes it in
# 1/6-character steps, so subtract the fixed 5/6 padding the host re-adds).
$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527 - 0.8333333333333333
$zhcn.Columns.Item(9).ColumnWidth = 40 - 0.8333333333333333
$zhcn.Columns.Item(10).ColumnWidth = 40 - 0.8333333333333333

# ---------------------------------------------------------------------------
# de-de sheet: same shape of change as zh-cn, but its own handback datetime.
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("I2").Value = "a0644925-81bd-4555-b74d-13847df737e1.md"
$dede.Range("J2").Value = "a0644925-81bd-4555-b74d-13847df737e1.01015c1d51e9dabb5a2674fb5c3e5b1e701e5437.de-de.xlf"
$dede.Range("K2").Value = "2016-08-29 00:25:18"

$dede.Range("I3").Value = "ed37b0c1-f71a-45ac-954a-94a9a70c5282.md"
$dede.Range("J3").Value = "ed37b0c1-f71a-45ac-954a-94a9a70c5282.d23391ddd3bd083b0c0eccc8d10aeaba46d43b77.de-de.xlf"
$dede.Range("K3").Value = "2016-08-29 00:25:18"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e9ddad79c30fb16c37f9428b0a1306f903cea0f6/e2e/a0644925-81bd-4555-b74d-13847df737e1.md", "", "", "a0644925-81bd-4555-b74d-13847df737e1.md")
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e9ddad79c30fb16c37f9428b0a1306f903cea0f6/e2e/a0644925-81bd-4555-b74d-13847df737e1.md", "", "", "a0644925-81bd-4555-b74d-13847df737e1.md")
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e9ddad79c30fb16c37f9428b0a1306f903cea0f6/e2e/ed37b0c1-f71a-45ac-954a-94a9a70c5282.md", "", "", "ed37b0c1-f71a-45ac-954a-94a9a70c5282.md")
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e9ddad79c30fb16c37f9428b0a1306f903cea0f6/e2e/ed37b0c1-f71a-45ac-954a-94a9a70c5282.md", "", "", "ed37b0c1-f71a-45ac-954a-94a9a70c5282.md")

$dede.Columns.Item(3).ColumnWidth = 29.9777047293527 - 0.8333333333333333
$dede.Columns.Item(9).ColumnWidth = 40 - 0.8333333333333333
$dede.Columns.Item(10).ColumnWidth = 40 - 0.8333333333333333

# ---------------------------------------------------------------------------
# Overview sheet column widths (zh-cn / de-de status columns got wider too)
# ---------------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527 - 0.8333333333333333
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527 - 0.8333333333333333

Write-Output "Generate Report for Handback: done"
